$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "1.00", "0.130")
# must be forced into Text so Excel stores the literal string instead of
# silently parsing/normalizing them into a number (dropping trailing zeros,
# switching to scientific notation, etc.) -- these columns hold
# pre-formatted display strings, not numeric data. The NumberFormat is
# restored to the default afterwards so the cell style is unaffected.

$ws.Range("D2").Value = "67.096.15"
$ws.Range("E2").Value = "  +4.40%  "
$ws.Range("D3").Value = "3.271.66"
$ws.Range("E3").Value = "  +2.97%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "3.269.58"
$ws.Range("E9").Value = "  +3.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.130"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.75"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("E12").Value = "  +5.03%  "
$ws.Range("D13").Value = "3.836.00"
$ws.Range("E13").Value = "  +3.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "67.061.08"
$ws.Range("E16").Value = "  +4.37%  "
$ws.Range("E17").Value = "  +3.37%  "
$ws.Range("D18").Value = "3.265.57"
$ws.Range("E18").Value = "  +3.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("D26").Value = "3.409.09"
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000120"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("E29").Value = "  +2.19%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.99"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.70"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("E35").Value = "  +6.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.85"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "167.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.51"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.859"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.59"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("D43").Value = "2.748.26"
$ws.Range("E43").Value = "  +4.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.46"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.03%  "
$ws.Range("E45").Value = "  +4.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "350.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.61%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0679"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.31%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0281"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.54%  "
$ws.Range("E51").Value = "  +2.39%  "
